# Weekly price-sheet update: a new week's Albahaca record for Vega Central
# Mapocho de Santiago is inserted as row 158, pushing the existing rows
# 158-224 down to 159-225 (dimension grows from A1:R224 to A1:R225).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 158; everything below shifts down one.
$ws.Rows("158:158").Insert()

# Populate the newly inserted row with this week's data point.
$ws.Range("A158").Value = 9
$ws.Range("B158").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C158").Value = "Metropolitana"
$ws.Range("D158").Value = 44510
$ws.Range("E158").Value = 13
$ws.Range("F158").Value = 100112052
$ws.Range("G158").Value = "Albahaca"
$ws.Range("H158").Value = "Sin especificar"
$ws.Range("I158").Value = "Primera"
$ws.Range("J158").Value = 61
$ws.Range("K158").Value = 5000
$ws.Range("L158").Value = 6000
$ws.Range("M158").Value = 5508
$ws.Range("N158").Value = "`$/docena de matas"
$ws.Range("O158").Value = "Provincia de Chacabuco"
$ws.Range("P158").Value = 918
$ws.Range("Q158").Value = 6
$ws.Range("R158").Value = "Hortaliza"
